# Append 12 new rows (225-236) to Sheet1, duplicating the content of row 224
# (question / model / response), except the very last row (236) whose
# response text differs slightly ("5" instead of "50" templates).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$question = "How many tadpole definitions am I allowed?"
$model = "llama3.2:latest"
$responseSame = "According to the document, you are allowed to create a maximum of 50 new tadpole definition templates."
$responseLast = "You are allowed to create a maximum of 5 new tadpole definition templates."

for ($r = 225; $r -le 236; $r++) {
    $ws.Cells.Item($r, 1).Value = $question
    $ws.Cells.Item($r, 2).Value = $model
    if ($r -eq 236) {
        $ws.Cells.Item($r, 3).Value = $responseLast
    } else {
        $ws.Cells.Item($r, 3).Value = $responseSame
    }
}
